$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.652559330171556
$ws.Range("C2").Value = 0.09191433375413993
$ws.Range("D2").Value = 0.04469045870509092
$ws.Range("F2").Value = 1.061829334936824
$ws.Range("G2").Value = 0.9210214308772464
$ws.Range("H2").Value = 0.9587367696032345
$ws.Range("I2").Value = 0.9516650563446944
$ws.Range("K2").Value = 0.4051525459359198
$ws.Range("L2").Value = 0.3030385400232376
$ws.Range("N2").Value = 1.803030238841753
$ws.Range("B3").Value = 0.6080022520366697
$ws.Range("C3").Value = 0.08968645335666281
$ws.Range("D3").Value = 0.04351573914959417
$ws.Range("F3").Value = 1.057662807273914
$ws.Range("G3").Value = 0.9184845056403219
$ws.Range("H3").Value = 0.9620002193509833
$ws.Range("I3").Value = 0.9548734544976369
$ws.Range("K3").Value = 0.3640989713066745
$ws.Range("L3").Value = 0.2918983375416104
$ws.Range("N3").Value = 1.821373653892056
$ws.Range("B4").Value = 0.5809078997976087
$ws.Range("C4").Value = 0.08829739446222362
$ws.Range("D4").Value = 0.04278456536061981
$ws.Range("F4").Value = 1.05570790097925
$ws.Range("G4").Value = 0.9174836469832428
$ws.Range("H4").Value = 0.9644185251876536
$ws.Range("I4").Value = 0.9573084924027739
$ws.Range("K4").Value = 0.3389880312932405
$ws.Range("L4").Value = 0.285232980978634
$ws.Range("N4").Value = 1.833227735542309
$ws.Range("B5").Value = 0.5699335643241739
$ws.Range("C5").Value = 0.08772602827294662
$ws.Range("D5").Value = 0.04248413311523436
$ws.Range("F5").Value = 1.055062964983662
$ws.Range("G5").Value = 0.9172157019943086
$ws.Range("H5").Value = 0.9655082954937058
$ws.Range("I5").Value = 0.9584177300971959
$ws.Range("K5").Value = 0.3287796261755318
$ws.Range("L5").Value = 0.2825607709324913
$ws.Range("N5").Value = 1.838206963934514
$ws.Range("B6").Value = 0.5681153349958663
$ws.Range("C6").Value = 0.08763083244492265
$ws.Range("D6").Value = 0.04243409762759498
$ws.Range("F6").Value = 1.054965035998862
$ws.Range("G6").Value = 0.9171796572922659
$ws.Range("H6").Value = 0.9656955515220318
$ws.Range("I6").Value = 0.9586089809262148
$ws.Range("K6").Value = 0.3270860183775142
$ws.Range("L6").Value = 0.2821197094509387
$ws.Range("N6").Value = 1.839042734128608
$ws.Range("B7").Value = 0.5807596248308471
$ws.Range("C7").Value = 0.08828971033268829
$ws.Range("D7").Value = 0.04278052362062823
$ws.Range("F7").Value = 1.055698588924152
$ws.Range("G7").Value = 0.9174794670212805
$ws.Range("H7").Value = 0.9644327998720286
$ws.Range("I7").Value = 0.9573229784932806
$ws.Range("K7").Value = 0.3388502574979668
$ws.Range("L7").Value = 0.2851967645194549
$ws.Range("N7").Value = 1.833294285583392
$ws.Range("B8").Value = 0.6371415554021098
$ws.Range("C8").Value = 0.09115054597480565
$ws.Range("D8").Value = 0.04428747654396403
$ws.Range("F8").Value = 1.060267475562156
$ws.Range("G8").Value = 0.9200310636994544
$ws.Range("H8").Value = 0.9597760008010852
$ws.Range("I8").Value = 0.952674782272787
$ws.Range("K8").Value = 0.3909775123810846
$ws.Range("L8").Value = 0.2991611372268892
$ws.Range("N8").Value = 1.809232323159662
$ws.Range("B9").Value = 0.749785083298093
$ws.Range("C9").Value = 0.09659318473334366
$ws.Range("D9").Value = 0.04716361999087582
$ws.Range("F9").Value = 1.07401712298271
$ws.Range("G9").Value = 0.9294592305258362
$ws.Range("H9").Value = 0.9539314907319323
$ws.Range("I9").Value = 0.9472508070648118
$ws.Range("K9").Value = 0.493952064491765
$ws.Range("L9").Value = 0.3279325064393532
$ws.Range("N9").Value = 1.766739193801136
$ws.Range("B10").Value = 0.8338003902358651
$ws.Range("C10").Value = 0.100490434492194
$ws.Range("D10").Value = 0.04922807156098941
$ws.Range("F10").Value = 1.087045777855607
$ws.Range("G10").Value = 0.9390945173581571
$ws.Range("H10").Value = 0.9516403366803274
$ws.Range("I10").Value = 0.9455185468309324
$ws.Range("K10").Value = 0.5700616198986665
$ws.Range("L10").Value = 0.3499199342520996
$ws.Range("N10").Value = 1.738381485768116
$ws.Range("B11").Value = 0.8722920706575792
$ws.Range("C11").Value = 0.1022415165524535
$ws.Range("D11").Value = 0.0501565854847712
$ws.Range("F11").Value = 1.093609991522385
$ws.Range("G11").Value = 0.9440685120818273
$ws.Range("H11").Value = 0.9510327284788502
$ws.Range("I11").Value = 0.9452202917692603
$ws.Range("K11").Value = 0.6047838554471525
$ws.Range("L11").Value = 0.3601078403928426
$ws.Range("N11").Value = 1.72610260430228
$ws.Range("B12").Value = 0.8869067239675701
$ws.Range("C12").Value = 0.1029014758922386
$ws.Range("D12").Value = 0.0505066509913803
$ws.Range("F12").Value = 1.096187422436714
$ws.Range("G12").Value = 0.9460371602859396
$ws.Range("H12").Value = 0.9508651163627491
$ws.Range("I12").Value = 0.9451778126477421
$ws.Range("K12").Value = 0.6179463562315561
$ws.Range("L12").Value = 0.3639924533920862
$ws.Range("N12").Value = 1.721542306045162
$ws.Range("B13").Value = 0.8837574845146037
$ws.Range("C13").Value = 0.1027594814295014
$ws.Range("D13").Value = 0.05043132698288844
$ws.Range("F13").Value = 1.095628247202285
$ws.Range("G13").Value = 0.9456093894011701
$ws.Range("H13").Value = 0.9508984364054101
$ws.Range("I13").Value = 0.9451838268070532
$ws.Range("K13").Value = 0.6151109620990667
$ws.Range("L13").Value = 0.3631546466507558
$ws.Range("N13").Value = 1.722520467678638
$ws.Range("B14").Value = 0.8734936545910728
$ws.Range("C14").Value = 0.1022958747007578
$ws.Range("D14").Value = 0.05018541659225662
$ws.Range("F14").Value = 1.093820200132001
$ws.Range("G14").Value = 0.9442287675339571
$ws.Range("H14").Value = 0.9510176869658125
$ws.Range("I14").Value = 0.9452153845782831
$ws.Range("K14").Value = 0.6058664650385595
$ws.Range("L14").Value = 0.3604268954371861
$ws.Range("N14").Value = 1.725725632444114
$ws.Range("B15").Value = 0.8672117853380996
$ws.Range("C15").Value = 0.1020114935134302
$ws.Range("D15").Value = 0.05003458801402161
$ws.Range("F15").Value = 1.092724663252156
$ws.Range("G15").Value = 0.9433941840326128
$ws.Range("H15").Value = 0.9510988668230738
$ws.Range("I15").Value = 0.9452438921661681
$ws.Range("K15").Value = 0.6002057482330656
$ws.Range("L15").Value = 0.358759541859726
$ws.Range("N15").Value = 1.727700540819402
$ws.Range("B16").Value = 0.831290310941796
$ws.Range("C16").Value = 0.1003755582947434
$ws.Range("D16").Value = 0.04916717600631415
$ws.Range("F16").Value = 1.086629621411831
$ws.Range("G16").Value = 0.9387813565952854
$ws.Range("H16").Value = 0.9516887893112624
$ws.Range("I16").Value = 0.945547896842541
$ws.Range("K16").Value = 0.5677944121301834
$ws.Range("L16").Value = 0.3492578650030822
$ws.Range("N16").Value = 1.739196449720872
$ws.Range("B17").Value = 0.8093231210298484
$ws.Range("C17").Value = 0.09936638021820698
$ws.Range("D17").Value = 0.04863231650600142
$ws.Range("F17").Value = 1.08305380058647
$ws.Range("G17").Value = 0.9361029681747084
$ws.Range("H17").Value = 0.9521619929569169
$ws.Range("I17").Value = 0.9458598589519696
$ws.Range("K17").Value = 0.5479363505317849
$ws.Range("L17").Value = 0.3434764378072259
$ws.Range("N17").Value = 1.746408003320827
$ws.Range("B18").Value = 0.7967138731609396
$ws.Range("C18").Value = 0.09878387628777574
$ws.Range("D18").Value = 0.04832368124693431
$ws.Range("F18").Value = 1.081057076607635
$ws.Range("G18").Value = 0.934618035631928
$ws.Range("H18").Value = 0.9524750799328672
$ws.Range("I18").Value = 0.9460853874844091
$ws.Range("K18").Value = 0.5365239179568846
$ws.Range("L18").Value = 0.3401685960926102
$ws.Range("N18").Value = 1.750614359296275
$ws.Range("B19").Value = 0.7924490282429986
$ws.Range("C19").Value = 0.098586298234693
$ws.Range("D19").Value = 0.04821901169723475
$ws.Range("F19").Value = 1.080391322717219
$ws.Range("G19").Value = 0.9341248090371437
$ws.Range("H19").Value = 0.952588113418571
$ws.Range("I19").Value = 0.946169663519747
$ws.Range("K19").Value = 0.5326614903111135
$ws.Range("L19").Value = 0.3390516212760275
$ws.Range("N19").Value = 1.752048596225169
$ws.Range("B20").Value = 0.8116589099676332
$ws.Range("C20").Value = 0.09947402128046434
$ws.Range("D20").Value = 0.04868935664100604
$ws.Range("F20").Value = 1.083428243368232
$ws.Range("G20").Value = 0.9363823312021538
$ws.Range("H20").Value = 0.9521073854930933
$ws.Range("I20").Value = 0.9458218790035531
$ws.Range("K20").Value = 0.5500493050641637
$ws.Range("L20").Value = 0.3440900714780355
$ws.Range("N20").Value = 1.74563426974392
$ws.Range("B21").Value = 0.8765073443983056
$ws.Range("C21").Value = 0.1024321324095183
$ws.Range("D21").Value = 0.05025768842718037
$ws.Range("F21").Value = 1.094348778150234
$ws.Range("G21").Value = 0.9446319788190465
$ws.Range("H21").Value = 0.9509809648117766
$ws.Range("I21").Value = 0.945204202638017
$ws.Range("K21").Value = 0.608581422524594
$ws.Range("L21").Value = 0.3612273784184197
$ws.Range("N21").Value = 1.724781769283826
$ws.Range("B22").Value = 0.9191148127430893
$ws.Range("C22").Value = 0.1043471402115301
$ws.Range("D22").Value = 0.05127368524825471
$ws.Range("F22").Value = 1.102020531064596
$ws.Range("G22").Value = 0.9505197186856691
$ws.Range("H22").Value = 0.9506089296493201
$ws.Range("I22").Value = 0.9452112486142781
$ws.Range("K22").Value = 0.6469166426923607
$ws.Range("L22").Value = 0.3725830609223095
$ws.Range("N22").Value = 1.711674883134918
$ws.Range("B23").Value = 0.8963539661709774
$ws.Range("C23").Value = 0.1033267389759658
$ws.Range("D23").Value = 0.05073225729785236
$ws.Range("F23").Value = 1.097877045325617
$ws.Range("G23").Value = 0.9473318799740582
$ws.Range("H23").Value = 0.9507741808958912
$ws.Range("I23").Value = 0.9451698936378889
$ws.Range("K23").Value = 0.6264491197625546
$ws.Range("L23").Value = 0.3665081052592427
$ws.Range("N23").Value = 1.718622523014417
$ws.Range("B24").Value = 0.8106028375077869
$ws.Range("C24").Value = 0.09942536396427926
$ws.Range("D24").Value = 0.04866357233712648
$ws.Range("F24").Value = 1.083258773810428
$ws.Range("G24").Value = 0.9362558601490036
$ws.Range("H24").Value = 0.9521319456977295
$ws.Range("I24").Value = 0.9458389058985475
$ws.Range("K24").Value = 0.5490940251125664
$ws.Range("L24").Value = 0.343812597926501
$ws.Range("N24").Value = 1.745983886747226
$ws.Range("B25").Value = 0.7190906345060171
$ws.Range("C25").Value = 0.09513866879973421
$ws.Range("D25").Value = 0.04639405485153958
$ws.Range("F25").Value = 1.069783983144859
$ws.Range("G25").Value = 0.9264338396792198
$ws.Range("H25").Value = 0.955160779497632
$ws.Range("I25").Value = 0.9483227168907433
$ws.Range("K25").Value = 0.4660146044783744
$ws.Range("L25").Value = 0.3200002199276355
$ws.Range("N25").Value = 1.777732230487075

Write-Output "done"
